# "added project idea to the list"
# Appends 8 new time-log rows (22-29) to the TODO sheet, covering
# 2024-10-15 and 2024-10-16. The Overview sheet recalculates on its own
# via existing SUMIFS formulas, so it needs no direct edits.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TODO")

# ---------------------------------------------------------------------
# 1) Clone per-column formatting from row 14 (a row that already has
#    every column B..J populated) down onto the new rows, but only onto
#    the specific cells that will actually hold data -- matching exactly
#    which columns are populated on each new row.
# ---------------------------------------------------------------------
$ws.Range("B14").Copy() | Out-Null
$ws.Range("B22:B29").PasteSpecial(-4122) | Out-Null

$ws.Range("C14").Copy() | Out-Null
$ws.Range("C23:C28").PasteSpecial(-4122) | Out-Null

$ws.Range("D14").Copy() | Out-Null
$ws.Range("D28").PasteSpecial(-4122) | Out-Null

$ws.Range("E14").Copy() | Out-Null
$ws.Range("E22:E28").PasteSpecial(-4122) | Out-Null

$ws.Range("F14").Copy() | Out-Null
$ws.Range("F22:F23").PasteSpecial(-4122) | Out-Null
$ws.Range("F25:F29").PasteSpecial(-4122) | Out-Null

$ws.Range("D14").Copy() | Out-Null
$ws.Range("F24").PasteSpecial(-4122) | Out-Null

$ws.Range("G14").Copy() | Out-Null
$ws.Range("G22:G29").PasteSpecial(-4122) | Out-Null

$ws.Range("H14").Copy() | Out-Null
$ws.Range("H22:H29").PasteSpecial(-4122) | Out-Null

$ws.Range("I14").Copy() | Out-Null
$ws.Range("I22:I29").PasteSpecial(-4122) | Out-Null

$ws.Range("J14").Copy() | Out-Null
$ws.Range("J22:J29").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2) Row heights: wrapped long descriptions push a couple of rows taller.
# ---------------------------------------------------------------------
$ws.Rows.Item(22).RowHeight = 34
$ws.Rows.Item(23).RowHeight = 17
$ws.Rows.Item(25).RowHeight = 17
$ws.Rows.Item(26).RowHeight = 34
$ws.Rows.Item(27).RowHeight = 17
$ws.Rows.Item(28).RowHeight = 17
$ws.Rows.Item(29).RowHeight = 17

# ---------------------------------------------------------------------
# 3) Data: dates, category/subcategory, task, description, status,
#    start/end time, and the elapsed-time formula per row.
# ---------------------------------------------------------------------

# Row 22 - Bredius rollers (continuation)
$ws.Range("B22").Value2 = 45580
$ws.Range("E22").Value2 = "Bredius rollers"
$ws.Range("F22").Value2 = "check voeg speler toe, maar server error 500, fix in code next time"
$ws.Range("G22").Value2 = "Completed"
$ws.Range("H22").Value2 = 0.64583333333333337
$ws.Range("I22").Value2 = 0.65625
$ws.Range("J22").Formula = "=IF(I22 < H22, I22 + 1, I22) - H22"

# Row 23 - Create message
$ws.Range("B23").Value2 = 45580
$ws.Range("C23").Value2 = "Univercity"
$ws.Range("E23").Value2 = "Create message"
$ws.Range("F23").Value2 = "Message to mohsen"
$ws.Range("G23").Value2 = "Completed"
$ws.Range("H23").Value2 = 0.65625
$ws.Range("I23").Value2 = 0.66666666666666663
$ws.Range("J23").Formula = "=IF(I23 < H23, I23 + 1, I23) - H23"

# Row 24 - Mac Github
$ws.Range("B24").Value2 = 45580
$ws.Range("C24").Value2 = "Project"
$ws.Range("E24").Value2 = "Mac Github"
$ws.Range("F24").Value2 = "Fix Personal Github MAC"
$ws.Range("G24").Value2 = "Completed"
$ws.Range("H24").Value2 = 0.66666666666666663
$ws.Range("I24").Value2 = 0.6875
$ws.Range("J24").Formula = "=IF(I24 < H24, I24 + 1, I24) - H24"

# Row 25 - Hockey
$ws.Range("B25").Value2 = 45580
$ws.Range("C25").Value2 = "Project"
$ws.Range("E25").Value2 = "Hockey"
$ws.Range("F25").Value2 = "Hockey Stick zoeken"
$ws.Range("G25").Value2 = "Completed"
$ws.Range("H25").Value2 = 0.6875
$ws.Range("I25").Value2 = 0.70833333333333337
$ws.Range("J25").Formula = "=IF(I25 < H25, I25 + 1, I25) - H25"

# Row 26 - QR Dive Mail
$ws.Range("B26").Value2 = 45580
$ws.Range("C26").Value2 = "SaaS"
$ws.Range("E26").Value2 = "QR Dive Mail"
$ws.Range("F26").Value2 = "QR dive mail fix end in spam | Only work for gmail not outlook"
$ws.Range("G26").Value2 = "Completed"
$ws.Range("H26").Value2 = 0.70833333333333337
$ws.Range("I26").Value2 = 0.8125
$ws.Range("J26").Formula = "=IF(I26 < H26, I26 + 1, I26) - H26"

# Row 27 - Xgino project
$ws.Range("B27").Value2 = 45581
$ws.Range("C27").Value2 = "Project"
$ws.Range("E27").Value2 = "Xgino project"
$ws.Range("F27").Value2 = "make project planning for xgino data science projects"
$ws.Range("G27").Value2 = "Completed"
$ws.Range("H27").Value2 = 0.41666666666666669
$ws.Range("I27").Value2 = 0.46875
$ws.Range("J27").Formula = "=IF(I27 < H27, I27 + 1, I27) - H27"

# Row 28 - QR Dive todo
$ws.Range("B28").Value2 = 45581
$ws.Range("C28").Value2 = "SaaS"
$ws.Range("D28").Value2 = "QR Dive"
$ws.Range("E28").Value2 = "QR Dive todo"
$ws.Range("F28").Value2 = "Clean up QR dive todo items and prioritize them"
$ws.Range("G28").Value2 = "Completed"
$ws.Range("H28").Value2 = 0.46875
$ws.Range("I28").Value2 = 0.48958333333333331
$ws.Range("J28").Formula = "=IF(I28 < H28, I28 + 1, I28) - H28"

# Row 29 - FOOD
$ws.Range("B29").Value2 = 45581
$ws.Range("F29").Value2 = "FOOD"
$ws.Range("G29").Value2 = "Completed"
$ws.Range("H29").Value2 = 0.48958333333333331
$ws.Range("I29").Value2 = 0.54166666666666663
$ws.Range("J29").Formula = "=IF(I29 < H29, I29 + 1, I29) - H29"

# ---------------------------------------------------------------------
# 4) Row 21's end-of-day entry shifted its start time earlier.
# ---------------------------------------------------------------------
$ws.Range("H21").Value2 = 0.5625

# ---------------------------------------------------------------------
# 5) Extend the "Status" (G) list-validation down to the new rows.
# ---------------------------------------------------------------------
$ws.Range("G13:G21").Validation.Delete()
$ws.Range("G13:G29").Validation.Add(3, 1, 1, "=Category!$D$3:$D$5")

# ---------------------------------------------------------------------
# 6) Selection / scroll position bookkeeping to match where the author
#    ended up after the edit.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("I30").Select()

$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("I8").Select()
